$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.62"
$ws.Range("E2").Value = "'3.54%"
$ws.Range("D3").Value = "'40.32"
$ws.Range("E3").Value = "'5.99%"
$ws.Range("D4").Value = "'5.811"
$ws.Range("E4").Value = "'12.30%"
$ws.Range("D5").Value = "'0.08012"
$ws.Range("E5").Value = "'0.33%"
$ws.Range("D6").Value = "'4.590"
$ws.Range("E6").Value = "'2.53%"
$ws.Range("D7").Value = "'8.717"
$ws.Range("E7").Value = "'2.31%"
$ws.Range("D8").Value = "'1.945"
$ws.Range("E8").Value = "'0.29%"
$ws.Range("D9").Value = "'2.941"
$ws.Range("E9").Value = "'-1.49%"
$ws.Range("D10").Value = "'0.9453"
$ws.Range("E10").Value = "'0.32%"
$ws.Range("D11").Value = "'0.1248"
$ws.Range("E11").Value = "'-0.38%"
$ws.Range("D12").Value = "'0.1962"
$ws.Range("E12").Value = "'1.04%"
$ws.Range("D13").Value = "'8.891"
$ws.Range("E13").Value = "'34.94%"
$ws.Range("D14").Value = "'0.09220"
$ws.Range("E14").Value = "'1.95%"
$ws.Range("D15").Value = "'0.03586"
$ws.Range("E15").Value = "'4.67%"
$ws.Range("D16").Value = "'0.09623"
$ws.Range("E16").Value = "'0.88%"
$ws.Range("D17").Value = "'0.001301"
$ws.Range("E17").Value = "'-4.36%"
$ws.Range("D18").Value = "'0.006199"
$ws.Range("E18").Value = "'2.35%"
$ws.Range("D19").Value = "'3.370"
$ws.Range("E19").Value = "'-1.53%"
$ws.Range("D20").Value = "'0.3528"
$ws.Range("E20").Value = "'0.37%"
$ws.Range("D21").Value = "'0.1407"
$ws.Range("E21").Value = "'7.79%"
$ws.Range("D22").Value = "'0.2422"
$ws.Range("E22").Value = "'5.15%"
$ws.Range("D23").Value = "'0.04404"
$ws.Range("E23").Value = "'0.77%"
$ws.Range("D24").Value = "'0.001263"
$ws.Range("E24").Value = "'3.04%"
$ws.Range("D25").Value = "'0.004314"
$ws.Range("E25").Value = "'-2.58%"
$ws.Range("E26").Value = "'-13.66%"
$ws.Range("E27").Value = "'0.41%"
$ws.Range("D39").Value = "'0.02425"
$ws.Range("E39").Value = "'-0.03%"
$ws.Range("D40").Value = "'0.05283"
$ws.Range("E40").Value = "'2.79%"
$ws.Range("D41").Value = "'0.007464"
$ws.Range("E41").Value = "'0.34%"
$ws.Range("D42").Value = "'0.1415"
$ws.Range("E42").Value = "'0.95%"
$ws.Range("D43").Value = "'0.008521"
$ws.Range("E43").Value = "'-0.91%"
$ws.Range("D44").Value = "'0.002108"
$ws.Range("E44").Value = "'3.69%"
$ws.Range("D45").Value = "'0.01086"
$ws.Range("E45").Value = "'24.07%"
$ws.Range("D46").Value = "'0.00006926"
$ws.Range("E46").Value = "'7.19%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.75%"
$ws.Range("D48").Value = "'0.003158"
$ws.Range("E48").Value = "'10.51%"
$ws.Range("D49").Value = "'0.001428"
$ws.Range("E49").Value = "'-15.19%"
$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.75%"
$ws.Range("D51").Value = "'0.0002008"
$ws.Range("E51").Value = "'0.75%"
